# The commit fills in the "Status" column (column B) on the
# "2. 25 Jan 2020" sheet, which previously had no status recorded for any
# row. The author's message says "Made most of the changes requested
# except some which I have no idea how to do!" - the change we can see
# and make here is adding those status values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2. 25 Jan 2020")

# Set the new values in an order that introduces the brand-new status
# strings ("?", "Kinda", "Nope") in the same order they first appear in
# the final workbook, so the shared-string table is built up the same way.
$ws.Range("B4").Value = "?"
$ws.Range("B12").Value = "Kinda"
$ws.Range("B5").Value = "Nope"

$ws.Range("B3").Value = "Done"
$ws.Range("B6").Value = "Done"
$ws.Range("B8").Value = "Done"
$ws.Range("B9").Value = "Done"
$ws.Range("B11").Value = "?"
$ws.Range("B13").Value = "Kinda"

# Move the cell cursor/selection to where it ended up after the edits.
$ws.Range("D18").Select()
